$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: was "Criação do Site por código" -> now "Script Banco de Dados"
# (done first so new shared-string "Script Banco de Dados" is registered
#  before "Feito" is introduced by the Row 9 edit below)
$ws.Range("C13").Value = "Script Banco de Dados"
$ws.Range("D13").Value = "Alta"
$ws.Range("E13").Value = "Feito"
$ws.Range("F13").Value = 45394
$ws.Range("G13").Value = "Larissa e Tabata"
$ws.Range("H13").Value = 1
$ws.Range("H13").NumberFormat = "0%"

# --- Row 7: Entregar novo protótipo de Calculadora -> set "Concluido %" to 50%
$ws.Range("H7").Value = 0.5
$ws.Range("H7").NumberFormat = "0%"

# --- Row 8: Modificações na documentação -> set "Concluido %" to 70%
$ws.Range("H8").Value = 0.7
$ws.Range("H8").NumberFormat = "0%"

# --- Row 9: Novo Banco de Dados com DER -> Status "Feito", Concluido % 100%
$ws.Range("E9").Value = "Feito"
$ws.Range("H9").Value = 1
$ws.Range("H9").NumberFormat = "0%"

# --- Row 10: Modificações nos Slides de Apresentação -> Concluido % 30%
$ws.Range("H10").Value = 0.3
$ws.Range("H10").NumberFormat = "0%"

# --- Row 11: Protótipo do Site no Figma -> Concluido % 0%
$ws.Range("H11").Value = 0
$ws.Range("H11").NumberFormat = "0%"

# --- Row 12: Verificação do código de arduino -> Concluido % 100%
$ws.Range("H12").Value = 1
$ws.Range("H12").NumberFormat = "0%"

# --- Update the active selection to H16
$ws.Range("H16").Select()

$wb.Save()
